{"js": "// Replace the date line and the 25 multiplication problems per the diff.\nconst replacements = [\n  [\"2024-11-18 Monday\", \"2024-11-19 Tuesday\"],\n  [\"93\u00d785=\", \"92\u00d784=\"],\n  [\"69\u00d715=\", \"87\u00d754=\"],\n  [\"31\u00d725=\", \"36\u00d753=\"],\n  [\"38\u00d750=\", \"32\u00d797=\"],\n  [\"99\u00d753=\", \"44\u00d712=\"],\n  [\"29\u00d786=\", \"71\u00d715=\"],\n  [\"27\u00d786=\", \"44\u00d724=\"],\n  [\"89\u00d717=\", \"96\u00d722=\"],\n  [\"33\u00d754=\", \"68\u00d739=\"],\n  [\"18\u00d732=\", \"16\u00d784=\"],\n  [\"29\u00d717=\", \"55\u00d744=\"],\n  [\"54\u00d758=\", \"19\u00d772=\"],\n  [\"64\u00d796=\", \"90\u00d793=\"],\n  [\"87\u00d746=\", \"55\u00d771=\"],\n  [\"34\u00d755=\", \"77\u00d765=\"],\n  [\"42\u00d712=\", \"84\u00d778=\"],\n  [\"46\u00d729=\", \"42\u00d779=\"],\n  [\"20\u00d742=\", \"21\u00d755=\"],\n  [\"58\u00d783=\", \"30\u00d727=\"],\n  [\"96\u00d798=\", \"51\u00d766=\"],\n  [\"69\u00d728=\", \"54\u00d715=\"],\n  [\"98\u00d721=\", \"35\u00d734=\"],\n  [\"21\u00d719=\", \"52\u00d715=\"],\n  [\"98\u00d752=\", \"42\u00d733=\"],\n  [\"91\u00d735=\", \"86\u00d736=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-11-18 Monday\", \"2024-11-19 Tuesday\"),\n    @(\"93\u00d785=\", \"92\u00d784=\"),\n    @(\"69\u00d715=\", \"87\u00d754=\"),\n    @(\"31\u00d725=\", \"36\u00d753=\"),\n    @(\"38\u00d750=\", \"32\u00d797=\"),\n    @(\"99\u00d753=\", \"44\u00d712=\"),\n    @(\"29\u00d786=\", \"71\u00d715=\"),\n    @(\"27\u00d786=\", \"44\u00d724=\"),\n    @(\"89\u00d717=\", \"96\u00d722=\"),\n    @(\"33\u00d754=\", \"68\u00d739=\"),\n    @(\"18\u00d732=\", \"16\u00d784=\"),\n    @(\"29\u00d717=\", \"55\u00d744=\"),\n    @(\"54\u00d758=\", \"19\u00d772=\"),\n    @(\"64\u00d796=\", \"90\u00d793=\"),\n    @(\"87\u00d746=\", \"55\u00d771=\"),\n    @(\"34\u00d755=\", \"77\u00d765=\"),\n    @(\"42\u00d712=\", \"84\u00d778=\"),\n    @(\"46\u00d729=\", \"42\u00d779=\"),\n    @(\"20\u00d742=\", \"21\u00d755=\"),\n    @(\"58\u00d783=\", \"30\u00d727=\"),\n    @(\"96\u00d798=\", \"51\u00d766=\"),\n    @(\"69\u00d728=\", \"54\u00d715=\"),\n    @(\"98\u00d721=\", \"35\u00d734=\"),\n    @(\"21\u00d719=\", \"52\u00d715=\"),\n    @(\"98\u00d752=\", \"42\u00d733=\"),\n    @(\"91\u00d735=\", \"86\u00d736=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute([ref]$old, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$new, [ref]2) | Out-Null\n}\n"}
